$wb = $excel.ActiveWorkbook
$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsData = $wb.Worksheets.Item("DATA")

# --- DATA sheet: insert a new "version" column between "browser" and "username" ---
$wsData.Columns.Item(4).Insert()

# Header
$wsData.Range("D1").Value = "version"

# Row 2 (loginLogoutTest / yes / chrome) now also executes against an explicit version
$wsData.Range("D2").Value = "79.0.3945.117"

# Row 3 used to be the "no" variant of loginLogoutTest; it becomes a second, executed
# variant of row 2 that targets a different chrome version
$wsData.Range("B3").Value = "yes"
$wsData.Range("D3").Value = "94.0.4606.61"

# Row 4 (failedLoginLogoutTest / yes / chrome) also gets an explicit version
$wsData.Range("D4").Value = "79.0.3945.117"

# Rows 5 and 6 are left without an explicit browser version (blank D5/D6)

# Match the quotePrefix/text style already used elsewhere in the workbook so the
# dotted version numbers are stored as plain text, same as RUNMANAGER!D2.
$wsRun.Range("D2").Copy()
$wsData.Range("D2:D4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Approximate the recalculated "best fit" column widths after the insert/shift.
$wsData.Columns.Item(3).ColumnWidth = 7
$wsData.Columns.Item(4).ColumnWidth = 12
$wsData.Columns.Item(5).ColumnWidth = 8
$wsData.Columns.Item(6).ColumnWidth = 9
$wsData.Columns.Item(7).ColumnWidth = 12.857142857142858

# Update the active selection to match the edited area
$wsData.Activate()
$wsData.Range("D5").Select()

Write-Output "done"
